$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 78
    3  = 80
    4  = 82
    5  = 83
    6  = 86
    7  = 89
    8  = 91
    9  = 93
    10 = 96
    11 = 98
    12 = 100
    13 = 103
    14 = 50
    15 = 75
    16 = 173
    17 = 194
    18 = 215
    19 = 243
    20 = 279
    21 = 300
    22 = 358
    23 = 370
    24 = 400
    25 = 417
    26 = 445
    27 = 509
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
